# 5.5.2.xlsx — add the "2020" data column (Q) to the indicator table.
#
# The sheet holds a small year-by-year table (years across row 4, values
# across row 5). A new year, 2020, with value 47.4, is appended in column Q,
# matching the formatting already used for the neighbouring 2019 column (P).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last existing year column (P) onto the new
# column (Q) before writing the new values, so Q4/Q5 pick up the same
# cell styles (right-aligned year header style, matching value style) as
# P4/P5 instead of a generic default style.
$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("P5").Copy($ws.Range("Q5"))

# New year column: header (2020) and its value (47.4).
$ws.Range("Q4").Value = 2020
$ws.Range("Q5").Value = 47.4

# Scroll the view over a bit and leave the selection on Q9, matching the
# saved view state of the uploaded workbook.
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("Q9").Select()
